$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new value would otherwise be
# auto-coerced to a number by the COM value setter (e.g. "1.00" -> 1),
# so they stay text just like the original inline-string cells.
$textCells = @("D4", "D5", "D6", "D8", "D10", "D12", "D14", "D17", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D43", "D44", "D45", "D47", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price / Volume(1h) figures (and the Bittensor/Cosmos
# row swap in B34:E35) from the latest GitHub Actions crypto refresh.
$ws.Range("D2").Value = "70.170.32"
$ws.Range("E2").Value = "  -3.21%  "
$ws.Range("D3").Value = "3.923.02"
$ws.Range("E3").Value = "  -2.96%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "533.26"
$ws.Range("E5").Value = "  +2.73%  "
$ws.Range("D6").Value = "146.35"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").Value = "3.916.91"
$ws.Range("E7").Value = "  -3.06%  "
$ws.Range("D8").Value = "0.681"
$ws.Range("E8").Value = "  -6.11%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "0.730"
$ws.Range("E11").Value = "  -5.87%  "
$ws.Range("D12").Value = "54.41"
$ws.Range("E12").Value = "  +15.51%  "
$ws.Range("E13").Value = "  -4.05%  "
$ws.Range("D14").Value = "10.46"
$ws.Range("E14").Value = "  -4.08%  "
$ws.Range("D15").Value = "4.545.26"
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("D16").Value = "3.922.46"
$ws.Range("E16").Value = "  -3.12%  "
$ws.Range("D17").Value = "20.23"
$ws.Range("E17").Value = "  -4.12%  "
$ws.Range("E18").Value = "  -3.71%  "
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("E20").Value = "  -4.30%  "
$ws.Range("D21").Value = "70.163.66"
$ws.Range("E21").Value = "  -3.01%  "
$ws.Range("D22").Value = "420.81"
$ws.Range("E22").Value = "  -5.44%  "
$ws.Range("D23").Value = "95.40"
$ws.Range("E23").Value = "  -8.97%  "
$ws.Range("D24").Value = "3.49"
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("D25").Value = "4.14"
$ws.Range("E25").Value = "  +3.60%  "
$ws.Range("D26").Value = "14.07"
$ws.Range("E26").Value = "  -3.78%  "
$ws.Range("D27").Value = "11.24"
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("D28").Value = "10.46"
$ws.Range("E28").Value = "  -5.11%  "
$ws.Range("D29").Value = "5.86"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("D30").Value = "3.63"
$ws.Range("E30").Value = "  +16.38%  "
$ws.Range("D31").Value = "36.04"
$ws.Range("E31").Value = "  -4.15%  "
$ws.Range("D32").Value = "7.60"
$ws.Range("E32").Value = "  +11.08%  "
$ws.Range("D33").Value = "48.99"
$ws.Range("E33").Value = "  +14.16%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "13.13"
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "675.46"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "0.127"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("D37").Value = "63.71"
$ws.Range("E37").Value = "  -5.66%  "
$ws.Range("D38").Value = "0.430"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "3.42"
$ws.Range("E39").Value = "  -5.58%  "
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").Value = "0.0₃0810"
$ws.Range("E41").Value = "  -5.81%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "3.19"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").Value = "0.0478"
$ws.Range("E45").Value = "  -3.56%  "
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("D47").Value = "0.147"
$ws.Range("E47").Value = "  -8.15%  "
$ws.Range("D48").Value = "9.56"
$ws.Range("E48").Value = "  +5.90%  "
$ws.Range("D49").Value = "3.33"
$ws.Range("E49").Value = "  -2.88%  "
$ws.Range("E50").Value = "  -4.17%  "
$ws.Range("E51").Value = "  +1.88%  "
